$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "292.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.45%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.90"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.898"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.18%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07289"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.70%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.287"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "24.76%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.680"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.08%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.720"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.15%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9006"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.59%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1666"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.60%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07902"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.77%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08023"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.28%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03097"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.74%"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.33%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001499"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.71%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005813"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.92%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.483"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.62%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.079"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.87%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3325"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.31%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.51%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.000"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-8.00%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2099"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.85%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04507"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.49%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001208"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.39%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "15.18%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001297"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.73%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01582"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.34%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04385"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.11%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007319"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.25%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1316"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.57%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001998"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.47%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009361"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-16.30%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005825"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.94%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.27%"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "4.29%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "20.55%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002095"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.27%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001996"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.27%"
